$d = $word.ActiveDocument

$replacements = @(
    @("2023-08-12 Saturday", "2023-08-13 Sunday"),
    @("67×29=1943", "67×52=3484"),
    @("79×41=3239", "25×97=2425"),
    @("64×73=4672", "84×93=7812"),
    @("14×43=602", "11×49=539"),
    @("15×95=1425", "48×79=3792"),
    @("28×34=952", "81×47=3807"),
    @("90×17=1530", "46×52=2392"),
    @("93×23=2139", "37×18=666"),
    @("48×66=3168", "43×11=473"),
    @("42×40=1680", "13×19=247"),
    @("78×23=1794", "91×40=3640"),
    @("51×50=2550", "80×74=5920"),
    @("46×13=598", "36×34=1224"),
    @("85×35=2975", "23×64=1472"),
    @("96×87=8352", "96×63=6048"),
    @("26×15=390", "97×60=5820"),
    @("95×17=1615", "33×84=2772"),
    @("84×64=5376", "71×81=5751"),
    @("31×19=589", "29×94=2726"),
    @("94×18=1692", "17×82=1394"),
    @("18×62=1116", "64×20=1280"),
    @("12×34=408", "37×45=1665"),
    @("12×30=360", "85×69=5865"),
    @("86×14=1204", "99×85=8415"),
    @("50×26=1300", "46×11=506")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
